# Generate Report for Handoff
# Appends two new file-status rows (9c32a983-... and f8d52f4d-...) to the
# Overview sheet and to each locale detail sheet (zh-cn, de-de), mirroring
# the layout already used for the existing 32f99451-... / c7dacf51-...
# rows: plain-text cell values (never real Excel dates, even though some
# look like dates), plus hyperlinks on the file-name / extension / target
# columns.

$wb = $excel.ActiveWorkbook

$newIds = @(
    "9c32a983-53c6-44dd-b9ce-ba1f2a0cd94c",
    "f8d52f4d-2026-4b61-b734-e4ffb3e10d19"
)
$newHashes = @(
    "9d99b061c50984289604842f1b13a1d45dc1d34e",
    "d3ff8fedcd88d540ce1bbd6174a4fe61c2ff124b"
)
$zhDatetimes = @("2016-03-15 03:17:00", "2016-03-15 03:17:00")
$deDatetimes = @("2016-03-15 03:17:08", "2016-03-15 03:17:08")
$handoffDate = "2016-17-15 03:17:08"
$status = "Ready for handoff"
$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/ab41e7463fa6c2543bd5ea66584f76f9fe2c72b0/e2e/"
$zhHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b87a00b893df17670bb96f11c2e4f9624da662f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/"
$deHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5f1f7f3de5fed9b82d719d86194b0e866276ea5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/"

function Set-CellStyle($range) {
    # Reuse the workbook's existing named "HyperLink" style so new linked
    # cells look like the pre-existing ones (blue + underline).
    $range.Style = "HyperLink"
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = 4 + $i
    $id = $newIds[$i]
    $mdName = "$id.md"

    $cellA = $wsOverview.Cells.Item($row, 1)
    $cellA.Value2 = $mdName
    $cellB = $wsOverview.Cells.Item($row, 2)
    $cellB.Value2 = $status
    $cellC = $wsOverview.Cells.Item($row, 3)
    $cellC.Value2 = $status
    $cellD = $wsOverview.Cells.Item($row, 4)
    $cellD.Value2 = $handoffDate

    $wsOverview.Hyperlinks.Add($cellA, ($mdBase + $mdName), "", "", $mdName) | Out-Null
    Set-CellStyle $cellA
}

# ---------------------------------------------------------------------
# Locale detail sheets: zh-cn / de-de
# Columns: A Source File Name | B File Extension | C Status |
#          D Latest Handoff File | E Latest Handoff Datetime |
#          F Latest Target File | G Latest Handback File |
#          H Latest Handback DateTime | I Handoff Reason |
#          J Dependency From | K Error Detail
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; HtBase = $zhHtBase; Datetimes = $zhDatetimes },
    @{ Sheet = "de-de"; Suffix = "de-de"; HtBase = $deHtBase; Datetimes = $deDatetimes }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    for ($i = 0; $i -lt $newIds.Length; $i++) {
        $row = 4 + $i
        $id = $newIds[$i]
        $hash = $newHashes[$i]
        $mdName = "$id.md"
        $xlfName = "$id.$hash.$($locale.Suffix).xlf"

        $cellA = $ws.Cells.Item($row, 1)
        $cellA.Value2 = $mdName
        $cellB = $ws.Cells.Item($row, 2)
        $cellB.Value2 = ".md"
        $cellC = $ws.Cells.Item($row, 3)
        $cellC.Value2 = $status
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.Value2 = $xlfName
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.Value2 = $locale.Datetimes[$i]
        $cellH = $ws.Cells.Item($row, 8)
        $cellH.Value2 = "0001-01-01 00:00:00"
        $cellI = $ws.Cells.Item($row, 9)
        $cellI.Value2 = "Include"

        $ws.Hyperlinks.Add($cellA, ($mdBase + $mdName), "", "", $mdName) | Out-Null
        Set-CellStyle $cellA

        $ws.Hyperlinks.Add($cellB, ($mdBase + $mdName), "", "", ".md") | Out-Null
        Set-CellStyle $cellB

        $ws.Hyperlinks.Add($cellD, ($locale.HtBase + $xlfName), "", "", $xlfName) | Out-Null
        Set-CellStyle $cellD
    }
}

Write-Host "Added handoff-ready rows for 9c32a983-... and f8d52f4d-... across Overview/zh-cn/de-de."
